$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 259; this shifts the existing rows
# 259-274 down to 260-275 (mirrors the diff, which shows every row from
# 259 onward taking on the values previously held by the row above it,
# with the former last row (274) duplicated into the new last row (275)).
$ws.Rows(259).Insert()

# Populate the newly inserted row 259 with its data.
$ws.Range("A259").Value = 8
$ws.Range("B259").Value = "Terminal La Palmera de La Serena"
$ws.Range("C259").Value = "Coquimbo"
$ws.Range("D259").Value = 44706
$ws.Range("E259").Value = 4
$ws.Range("F259").Value = 100112003
$ws.Range("G259").Value = "Ajo"
$ws.Range("H259").Value = "Chino"
$ws.Range("I259").Value = "Primera"
$ws.Range("J259").Value = 440
$ws.Range("K259").Value = 18500
$ws.Range("L259").Value = 19000
$ws.Range("M259").Value = 18750
$ws.Range("N259").Value = "$/caja 10 kilos"
$ws.Range("O259").Value = "China"
$ws.Range("P259").Value = 1875
$ws.Range("Q259").Value = 10
$ws.Range("R259").Value = "Hortaliza"
